$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''46.379.21'
$ws.Range("E2").Value = '  +2.04%  '
$ws.Range("D3").Value = '''2.597.82'
$ws.Range("E3").Value = '  +7.07%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '''305.50'
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("D6").Value = '''99.62'
$ws.Range("E6").Value = '  +5.25%  '
$ws.Range("D7").Value = '''0.600'
$ws.Range("E7").Value = '  +6.09%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '''0.581'
$ws.Range("E9").Value = '  +15.52%  '
$ws.Range("D10").Value = '''39.22'
$ws.Range("E10").Value = '  +12.56%  '
$ws.Range("D11").Value = '''54.42'
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("D12").Value = '''0.0841'
$ws.Range("E12").Value = '  +7.73%  '
$ws.Range("D13").Value = '''8.19'
$ws.Range("E13").Value = '  +16.02%  '
$ws.Range("D14").Value = '''2.990.24'
$ws.Range("E14").Value = '  +6.87%  '
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").Value = '''2.652.20'
$ws.Range("E16").Value = '  +9.00%  '
$ws.Range("D17").Value = '''0.914'
$ws.Range("E17").Value = '  +8.75%  '
$ws.Range("D18").Value = '''14.91'
$ws.Range("E18").Value = '  +5.58%  '
$ws.Range("D19").Value = '''46.452.41'
$ws.Range("E19").Value = '  +2.33%  '
$ws.Range("D20").Value = '''0.0000101'
$ws.Range("E20").Value = '  +7.17%  '
$ws.Range("D21").Value = '''12.97'
$ws.Range("E21").Value = '  +4.68%  '
$ws.Range("D22").Value = '''6.66'
$ws.Range("E22").Value = '  +7.99%  '
$ws.Range("D23").Value = '''71.98'
$ws.Range("E23").Value = '  +7.16%  '
$ws.Range("D24").Value = '''272.39'
$ws.Range("E24").Value = '  +13.09%  '
$ws.Range("D25").Value = '''3.02'
$ws.Range("E25").Value = '  +8.17%  '
$ws.Range("D26").Value = '''30.28'
$ws.Range("E26").Value = '  +42.53%  '
$ws.Range("D27").Value = '''2.16'
$ws.Range("E27").Value = '  +11.41%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").Value = '''10.53'
$ws.Range("E30").Value = '  +8.87%  '
$ws.Range("D31").Value = '''2.31'
$ws.Range("E31").Value = '  +3.57%  '
$ws.Range("D32").Value = '''39.07'
$ws.Range("E32").Value = '  +2.12%  '
$ws.Range("D33").Value = '''6.18'
$ws.Range("E33").Value = '  +12.97%  '
$ws.Range("D34").Value = '''3.64'
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("E35").Value = '  +4.22%  '
$ws.Range("D36").Value = '''0.0839'
$ws.Range("E36").Value = '  +9.77%  '
$ws.Range("D37").Value = '''2.18'
$ws.Range("E37").Value = '  +9.71%  '
$ws.Range("D38").Value = '''149.67'
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").Value = '''0.123'
$ws.Range("E39").Value = '  +8.14%  '
$ws.Range("D40").Value = '''0.122'
$ws.Range("E40").Value = '  +6.11%  '
$ws.Range("D41").Value = '''23.20'
$ws.Range("E41").Value = '  +43.47%  '
$ws.Range("D42").Value = '''16.05'
$ws.Range("E42").Value = '  +8.67%  '
$ws.Range("D43").Value = '''0.0329'
$ws.Range("E43").Value = '  +10.52%  '
$ws.Range("D44").Value = '''4.09'
$ws.Range("E44").Value = '  +7.84%  '
$ws.Range("D45").Value = '''3.59'
$ws.Range("E45").Value = '  +11.85%  '
$ws.Range("D46").Value = '''2.135.87'
$ws.Range("E46").Value = '  +5.91%  '
$ws.Range("D47").Value = '''0.997'
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").Value = '''93.99'
$ws.Range("E48").Value = '  +5.68%  '
$ws.Range("D49").Value = '''9.67'
$ws.Range("E49").Value = '  +12.46%  '
$ws.Range("D50").Value = '''1.78'
$ws.Range("E50").Value = '  +3.06%  '
$ws.Range("D51").Value = '''108.67'
$ws.Range("E51").Value = '  +7.63%  '
